# EPSSprep.xlsx maintenance edit:
#  - Strip the stray leading spaces that had crept into a couple of the
#    CVE id cells (B12, B16) so they match the rest of the column and
#    de-duplicate against the existing "clean" shared-string entries.
#  - Strip the stray leading (non-breaking) space in the CVSS vector
#    cell for CVE-2022-29464 (C14) for the same reason.
#  - Leave the cursor/selection parked past the data (B20), matching
#    where the author left off after the cleanup pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B12").Value = "CVE-2021-27877"
$ws.Range("B16").Value = "CVE-2023-27532"
$ws.Range("C14").Value = "CVSS:3.1/AV:N/AC:L/PR:N/UI:N/S:U/C:H/I:H/A:H"

$ws.Range("B20").Select() | Out-Null
